# Auto-generated edit script for BRVM recommandations workbook update
# 🔄 MAJ automatique BRVM via GitHub Actions
# Refreshes both the "Recommandations" leaderboard (sorted by Variation Totale)
# and the "Top_YTD" leaderboard (sorted by Progression YTD) with the latest data.

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd  = $wb.Worksheets.Item("Top_YTD")

# --- Recommandations sheet (columns A:G, rows 2-39) ---
# Row 2: BRVM - SERVICES PUBLICS
$wsReco.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$wsReco.Cells.Item(2, 2).Value = 0
$wsReco.Cells.Item(2, 3).Value = 8
$wsReco.Cells.Item(2, 4).Value = 3361.5
$wsReco.Cells.Item(2, 5).Value = 115.65
$wsReco.Cells.Item(2, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(2, 7).Value = "➖ Neutre"

# Row 3: NEI-CEDA CI
$wsReco.Cells.Item(3, 1).Value = "NEI-CEDA CI"
$wsReco.Cells.Item(3, 2).Value = 0
$wsReco.Cells.Item(3, 3).Value = 3
$wsReco.Cells.Item(3, 4).Value = 2910
$wsReco.Cells.Item(3, 5).Value = 995
$wsReco.Cells.Item(3, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(3, 7).Value = "➖ Neutre"

# Row 4: AIR LIQUIDE CI
$wsReco.Cells.Item(4, 1).Value = "AIR LIQUIDE CI"
$wsReco.Cells.Item(4, 2).Value = 0
$wsReco.Cells.Item(4, 3).Value = 4
$wsReco.Cells.Item(4, 4).Value = 2800
$wsReco.Cells.Item(4, 5).Value = 695
$wsReco.Cells.Item(4, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(4, 7).Value = "➖ Neutre"

# Row 5: BRVM - AUTRES SECTEURS
$wsReco.Cells.Item(5, 1).Value = "BRVM - AUTRES SECTEURS"
$wsReco.Cells.Item(5, 2).Value = 0
$wsReco.Cells.Item(5, 3).Value = 4
$wsReco.Cells.Item(5, 4).Value = 2429.89
$wsReco.Cells.Item(5, 5).Value = 601.76
$wsReco.Cells.Item(5, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(5, 7).Value = "➖ Neutre"

# Row 6: BRVM - DISTRIBUTION
$wsReco.Cells.Item(6, 1).Value = "BRVM - DISTRIBUTION"
$wsReco.Cells.Item(6, 2).Value = 0
$wsReco.Cells.Item(6, 3).Value = 4
$wsReco.Cells.Item(6, 4).Value = 2008.6
$wsReco.Cells.Item(6, 5).Value = 498.22
$wsReco.Cells.Item(6, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(6, 7).Value = "➖ Neutre"

# Row 7: BRVM - TRANSPORT
$wsReco.Cells.Item(7, 1).Value = "BRVM - TRANSPORT"
$wsReco.Cells.Item(7, 2).Value = 0
$wsReco.Cells.Item(7, 3).Value = 4
$wsReco.Cells.Item(7, 4).Value = 1460.07
$wsReco.Cells.Item(7, 5).Value = 362.26
$wsReco.Cells.Item(7, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(7, 7).Value = "➖ Neutre"

# Row 8: BRVM - AGRICULTURE
$wsReco.Cells.Item(8, 1).Value = "BRVM - AGRICULTURE"
$wsReco.Cells.Item(8, 2).Value = 0
$wsReco.Cells.Item(8, 3).Value = 4
$wsReco.Cells.Item(8, 4).Value = 1334.31
$wsReco.Cells.Item(8, 5).Value = 332.49
$wsReco.Cells.Item(8, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(8, 7).Value = "➖ Neutre"

# Row 9: BRVM - CONSOMMATION DISCRETIONNAIRE
$wsReco.Cells.Item(9, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsReco.Cells.Item(9, 2).Value = 0
$wsReco.Cells.Item(9, 3).Value = 4
$wsReco.Cells.Item(9, 4).Value = 711.53
$wsReco.Cells.Item(9, 5).Value = 173.91
$wsReco.Cells.Item(9, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(9, 7).Value = "➖ Neutre"

# Row 10: BRVM - FINANCES
$wsReco.Cells.Item(10, 1).Value = "BRVM - FINANCES"
$wsReco.Cells.Item(10, 2).Value = 0
$wsReco.Cells.Item(10, 3).Value = 4
$wsReco.Cells.Item(10, 4).Value = 573.95
$wsReco.Cells.Item(10, 5).Value = 144.87
$wsReco.Cells.Item(10, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(10, 7).Value = "➖ Neutre"

# Row 11: BRVM - SERVICES FINANCIERS
$wsReco.Cells.Item(11, 1).Value = "BRVM - SERVICES FINANCIERS"
$wsReco.Cells.Item(11, 2).Value = 0
$wsReco.Cells.Item(11, 3).Value = 4
$wsReco.Cells.Item(11, 4).Value = 564.07
$wsReco.Cells.Item(11, 5).Value = 142.37
$wsReco.Cells.Item(11, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(11, 7).Value = "➖ Neutre"

# Row 12: BRVM-PRESTIGE
$wsReco.Cells.Item(12, 1).Value = "BRVM-PRESTIGE"
$wsReco.Cells.Item(12, 2).Value = 0
$wsReco.Cells.Item(12, 3).Value = 4
$wsReco.Cells.Item(12, 4).Value = 563.87
$wsReco.Cells.Item(12, 5).Value = 142.87
$wsReco.Cells.Item(12, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(12, 7).Value = "➖ Neutre"

# Row 13: BRVM - INDUSTRIELS
$wsReco.Cells.Item(13, 1).Value = "BRVM - INDUSTRIELS"
$wsReco.Cells.Item(13, 2).Value = 0
$wsReco.Cells.Item(13, 3).Value = 4
$wsReco.Cells.Item(13, 4).Value = 505.22
$wsReco.Cells.Item(13, 5).Value = 124.51
$wsReco.Cells.Item(13, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(13, 7).Value = "➖ Neutre"

# Row 14: BRVM - ENERGIE
$wsReco.Cells.Item(14, 1).Value = "BRVM - ENERGIE"
$wsReco.Cells.Item(14, 2).Value = 0
$wsReco.Cells.Item(14, 3).Value = 4
$wsReco.Cells.Item(14, 4).Value = 432.67
$wsReco.Cells.Item(14, 5).Value = 109.52
$wsReco.Cells.Item(14, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(14, 7).Value = "➖ Neutre"

# Row 15: BRVM - TELECOMMUNICATIONS
$wsReco.Cells.Item(15, 1).Value = "BRVM - TELECOMMUNICATIONS"
$wsReco.Cells.Item(15, 2).Value = 0
$wsReco.Cells.Item(15, 3).Value = 4
$wsReco.Cells.Item(15, 4).Value = 377.31
$wsReco.Cells.Item(15, 5).Value = 94.29
$wsReco.Cells.Item(15, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(15, 7).Value = "➖ Neutre"

# Row 16: BRVM - INDUSTRIE                (**)
$wsReco.Cells.Item(16, 1).Value = "BRVM - INDUSTRIE                (**)"
$wsReco.Cells.Item(16, 2).Value = 0
$wsReco.Cells.Item(16, 3).Value = 1
$wsReco.Cells.Item(16, 4).Value = 235.63
$wsReco.Cells.Item(16, 5).Value = 235.63
$wsReco.Cells.Item(16, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(16, 7).Value = "➖ Neutre"

# Row 17: BRVM-PRINCIPAL                   (**)
$wsReco.Cells.Item(17, 1).Value = "BRVM-PRINCIPAL                   (**)"
$wsReco.Cells.Item(17, 2).Value = 0
$wsReco.Cells.Item(17, 3).Value = 1
$wsReco.Cells.Item(17, 4).Value = 210.45
$wsReco.Cells.Item(17, 5).Value = 210.45
$wsReco.Cells.Item(17, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(17, 7).Value = "➖ Neutre"

# Row 18: BRVM - CONSOMMATION DE BASE         (**)
$wsReco.Cells.Item(18, 1).Value = "BRVM - CONSOMMATION DE BASE         (**)"
$wsReco.Cells.Item(18, 2).Value = 0
$wsReco.Cells.Item(18, 3).Value = 1
$wsReco.Cells.Item(18, 4).Value = 200.13
$wsReco.Cells.Item(18, 5).Value = 200.13
$wsReco.Cells.Item(18, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(18, 7).Value = "➖ Neutre"

# Row 19: UNILEVER CI (UNLC)
$wsReco.Cells.Item(19, 1).Value = "UNILEVER CI (UNLC)"
$wsReco.Cells.Item(19, 2).Value = 4
$wsReco.Cells.Item(19, 3).Value = 0
$wsReco.Cells.Item(19, 4).Value = 29.96
$wsReco.Cells.Item(19, 5).Value = 7.49
$wsReco.Cells.Item(19, 6).Value = "🟢 Achat"
$wsReco.Cells.Item(19, 7).Value = "✅ Renforcer"

# Row 20: NESTLE CI (NTLC)
$wsReco.Cells.Item(20, 1).Value = "NESTLE CI (NTLC)"
$wsReco.Cells.Item(20, 2).Value = 2
$wsReco.Cells.Item(20, 3).Value = 0
$wsReco.Cells.Item(20, 4).Value = 8.93
$wsReco.Cells.Item(20, 5).Value = 4.58
$wsReco.Cells.Item(20, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(20, 7).Value = "➖ Neutre"

# Row 21: SAFCA CI (SAFC)
$wsReco.Cells.Item(21, 1).Value = "SAFCA CI (SAFC)"
$wsReco.Cells.Item(21, 2).Value = 2
$wsReco.Cells.Item(21, 3).Value = 2
$wsReco.Cells.Item(21, 4).Value = 8.65
$wsReco.Cells.Item(21, 5).Value = -2.07
$wsReco.Cells.Item(21, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(21, 7).Value = "👀 À surveiller"

# Row 22: SICABLE CI (CABC)
$wsReco.Cells.Item(22, 1).Value = "SICABLE CI (CABC)"
$wsReco.Cells.Item(22, 2).Value = 2
$wsReco.Cells.Item(22, 3).Value = 1
$wsReco.Cells.Item(22, 4).Value = 7.45
$wsReco.Cells.Item(22, 5).Value = -7.45
$wsReco.Cells.Item(22, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(22, 7).Value = "👀 À surveiller"

# Row 23: SETAO CI (STAC)
$wsReco.Cells.Item(23, 1).Value = "SETAO CI (STAC)"
$wsReco.Cells.Item(23, 2).Value = 2
$wsReco.Cells.Item(23, 3).Value = 1
$wsReco.Cells.Item(23, 4).Value = 6.69
$wsReco.Cells.Item(23, 5).Value = -4.17
$wsReco.Cells.Item(23, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(23, 7).Value = "👀 À surveiller"

# Row 24: BERNABE CI (BNBC)
$wsReco.Cells.Item(24, 1).Value = "BERNABE CI (BNBC)"
$wsReco.Cells.Item(24, 2).Value = 1
$wsReco.Cells.Item(24, 3).Value = 0
$wsReco.Cells.Item(24, 4).Value = 5.32
$wsReco.Cells.Item(24, 5).Value = 5.32
$wsReco.Cells.Item(24, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(24, 7).Value = "➖ Neutre"

# Row 25: TRACTAFRIC MOTORS CI (PRSC)
$wsReco.Cells.Item(25, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$wsReco.Cells.Item(25, 2).Value = 1
$wsReco.Cells.Item(25, 3).Value = 1
$wsReco.Cells.Item(25, 4).Value = 4.7
$wsReco.Cells.Item(25, 5).Value = -2.73
$wsReco.Cells.Item(25, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(25, 7).Value = "👀 À surveiller"

# Row 26: VIVO ENERGY CI (SHEC)
$wsReco.Cells.Item(26, 1).Value = "VIVO ENERGY CI (SHEC)"
$wsReco.Cells.Item(26, 2).Value = 1
$wsReco.Cells.Item(26, 3).Value = 0
$wsReco.Cells.Item(26, 4).Value = 4.45
$wsReco.Cells.Item(26, 5).Value = 4.45
$wsReco.Cells.Item(26, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(26, 7).Value = "➖ Neutre"

# Row 27: SICOR CI (SICC)
$wsReco.Cells.Item(27, 1).Value = "SICOR CI (SICC)"
$wsReco.Cells.Item(27, 2).Value = 1
$wsReco.Cells.Item(27, 3).Value = 1
$wsReco.Cells.Item(27, 4).Value = 2.05
$wsReco.Cells.Item(27, 5).Value = 7.45
$wsReco.Cells.Item(27, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(27, 7).Value = "👀 À surveiller"

# Row 28: ECOBANK TRANS. INCORP. TG (ETIT)
$wsReco.Cells.Item(28, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$wsReco.Cells.Item(28, 2).Value = 1
$wsReco.Cells.Item(28, 3).Value = 1
$wsReco.Cells.Item(28, 4).Value = 0.2
$wsReco.Cells.Item(28, 5).Value = 4.55
$wsReco.Cells.Item(28, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(28, 7).Value = "👀 À surveiller"

# Row 29: ONATEL BF (ONTBF)
$wsReco.Cells.Item(29, 1).Value = "ONATEL BF (ONTBF)"
$wsReco.Cells.Item(29, 2).Value = 1
$wsReco.Cells.Item(29, 3).Value = 1
$wsReco.Cells.Item(29, 4).Value = 0.15
$wsReco.Cells.Item(29, 5).Value = 3.96
$wsReco.Cells.Item(29, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(29, 7).Value = "👀 À surveiller"

# Row 30: TOTAL
$wsReco.Cells.Item(30, 1).Value = "TOTAL"
$wsReco.Cells.Item(30, 2).Value = 0
$wsReco.Cells.Item(30, 3).Value = 3
$wsReco.Cells.Item(30, 4).Value = 0
$wsReco.Cells.Item(30, 5).Value = 0
$wsReco.Cells.Item(30, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(30, 7).Value = "➖ Neutre"

# Row 31: BICI CI (BICC)
$wsReco.Cells.Item(31, 1).Value = "BICI CI (BICC)"
$wsReco.Cells.Item(31, 2).Value = 1
$wsReco.Cells.Item(31, 3).Value = 1
$wsReco.Cells.Item(31, 4).Value = -0.45
$wsReco.Cells.Item(31, 5).Value = 2.05
$wsReco.Cells.Item(31, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(31, 7).Value = "👀 À surveiller"

# Row 32: SOLIBRA CI (SLBC)
$wsReco.Cells.Item(32, 1).Value = "SOLIBRA CI (SLBC)"
$wsReco.Cells.Item(32, 2).Value = 0
$wsReco.Cells.Item(32, 3).Value = 1
$wsReco.Cells.Item(32, 4).Value = -1.2
$wsReco.Cells.Item(32, 5).Value = -1.2
$wsReco.Cells.Item(32, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(32, 7).Value = "➖ Neutre"

# Row 33: SOGB CI (SOGC)
$wsReco.Cells.Item(33, 1).Value = "SOGB CI (SOGC)"
$wsReco.Cells.Item(33, 2).Value = 0
$wsReco.Cells.Item(33, 3).Value = 1
$wsReco.Cells.Item(33, 4).Value = -1.94
$wsReco.Cells.Item(33, 5).Value = -1.94
$wsReco.Cells.Item(33, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(33, 7).Value = "➖ Neutre"

# Row 34: SUCRIVOIRE (SCRC)
$wsReco.Cells.Item(34, 1).Value = "SUCRIVOIRE (SCRC)"
$wsReco.Cells.Item(34, 2).Value = 0
$wsReco.Cells.Item(34, 3).Value = 1
$wsReco.Cells.Item(34, 4).Value = -3.21
$wsReco.Cells.Item(34, 5).Value = -3.21
$wsReco.Cells.Item(34, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(34, 7).Value = "➖ Neutre"

# Row 35: CIE CI (CIEC)
$wsReco.Cells.Item(35, 1).Value = "CIE CI (CIEC)"
$wsReco.Cells.Item(35, 2).Value = 0
$wsReco.Cells.Item(35, 3).Value = 1
$wsReco.Cells.Item(35, 4).Value = -3.27
$wsReco.Cells.Item(35, 5).Value = -3.27
$wsReco.Cells.Item(35, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(35, 7).Value = "➖ Neutre"

# Row 36: CFAO MOTORS CI (CFAC)
$wsReco.Cells.Item(36, 1).Value = "CFAO MOTORS CI (CFAC)"
$wsReco.Cells.Item(36, 2).Value = 0
$wsReco.Cells.Item(36, 3).Value = 1
$wsReco.Cells.Item(36, 4).Value = -4.97
$wsReco.Cells.Item(36, 5).Value = -4.97
$wsReco.Cells.Item(36, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(36, 7).Value = "➖ Neutre"

# Row 37: FILTISAC CI (FTSC)
$wsReco.Cells.Item(37, 1).Value = "FILTISAC CI (FTSC)"
$wsReco.Cells.Item(37, 2).Value = 0
$wsReco.Cells.Item(37, 3).Value = 2
$wsReco.Cells.Item(37, 4).Value = -6.09
$wsReco.Cells.Item(37, 5).Value = -2.89
$wsReco.Cells.Item(37, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(37, 7).Value = "➖ Neutre"

# Row 38: ORAGROUP TOGO (ORGT)
$wsReco.Cells.Item(38, 1).Value = "ORAGROUP TOGO (ORGT)"
$wsReco.Cells.Item(38, 2).Value = 0
$wsReco.Cells.Item(38, 3).Value = 2
$wsReco.Cells.Item(38, 4).Value = -6.79
$wsReco.Cells.Item(38, 5).Value = -4.12
$wsReco.Cells.Item(38, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(38, 7).Value = "➖ Neutre"

# Row 39: AFRICA GLOBAL LOGISTICS CI (SDSC)
$wsReco.Cells.Item(39, 1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$wsReco.Cells.Item(39, 2).Value = 0
$wsReco.Cells.Item(39, 3).Value = 2
$wsReco.Cells.Item(39, 4).Value = -14.65
$wsReco.Cells.Item(39, 5).Value = -7.3
$wsReco.Cells.Item(39, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(39, 7).Value = "➖ Neutre"

# --- Top_YTD sheet (columns A:B, rows 2-11) ---
# Row 2: BRVM - SERVICES PUBLICS
$wsYtd.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$wsYtd.Cells.Item(2, 2).Value = 9872644.53

# Row 3: AIR LIQUIDE CI
$wsYtd.Cells.Item(3, 1).Value = "AIR LIQUIDE CI"
$wsYtd.Cells.Item(3, 2).Value = 409484

# Row 4: BRVM - AUTRES SECTEURS
$wsYtd.Cells.Item(4, 1).Value = "BRVM - AUTRES SECTEURS"
$wsYtd.Cells.Item(4, 2).Value = 250384.55

# Row 5: BRVM - DISTRIBUTION
$wsYtd.Cells.Item(5, 1).Value = "BRVM - DISTRIBUTION"
$wsYtd.Cells.Item(5, 2).Value = 131343.15

# Row 6: NEI-CEDA CI
$wsYtd.Cells.Item(6, 1).Value = "NEI-CEDA CI"
$wsYtd.Cells.Item(6, 2).Value = 122271.72

# Row 7: BRVM - TRANSPORT
$wsYtd.Cells.Item(7, 1).Value = "BRVM - TRANSPORT"
$wsYtd.Cells.Item(7, 2).Value = 46597.77

# Row 8: BRVM - AGRICULTURE
$wsYtd.Cells.Item(8, 1).Value = "BRVM - AGRICULTURE"
$wsYtd.Cells.Item(8, 2).Value = 35239.61

# Row 9: BRVM - CONSOMMATION DISCRETIONNAIRE
$wsYtd.Cells.Item(9, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsYtd.Cells.Item(9, 2).Value = 5860.13

# Row 10: BRVM - FINANCES
$wsYtd.Cells.Item(10, 1).Value = "BRVM - FINANCES"
$wsYtd.Cells.Item(10, 2).Value = 3414.77

# Row 11: BRVM - SERVICES FINANCIERS
$wsYtd.Cells.Item(11, 1).Value = "BRVM - SERVICES FINANCIERS"
$wsYtd.Cells.Item(11, 2).Value = 3274.31

